$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Cells.Item(76, 8).Value = 3997.5
$ws.Cells.Item(76, 9).Value = 3000
$ws.Cells.Item(76, 10).Value = 4995
$ws.Cells.Item(76, 11).Value = 3000
$ws.Cells.Item(76, 12).Value = 4995
$ws.Cells.Item(76, 13).Value = -2685
$ws.Cells.Item(76, 14).Value = -5625
# Row 79
$ws.Cells.Item(79, 8).Value = 3997.5
$ws.Cells.Item(79, 9).Value = 3000
$ws.Cells.Item(79, 10).Value = 4995
$ws.Cells.Item(79, 11).Value = 3000
$ws.Cells.Item(79, 12).Value = 4995
$ws.Cells.Item(79, 13).Value = -1908
$ws.Cells.Item(79, 14).Value = -7179
# Row 135
$ws.Cells.Item(135, 8).Value = 961.93335
$ws.Cells.Item(135, 9).Value = 616.7727
$ws.Cells.Item(135, 11).Value = 5550.954299999999
$ws.Cells.Item(135, 13).Value = -3015.954299999999
# Row 137
$ws.Cells.Item(137, 8).Value = 2166.6667
$ws.Cells.Item(137, 9).Value = 2000
$ws.Cells.Item(137, 10).Value = 2250
$ws.Cells.Item(137, 11).Value = 6000
$ws.Cells.Item(137, 12).Value = 6750
$ws.Cells.Item(137, 13).Value = -3450
$ws.Cells.Item(137, 14).Value = -11850

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 2985.3381
$ws.Cells.Item(32, 9).Value = 2216.242
$ws.Cells.Item(32, 11).Value = 2216.242
$ws.Cells.Item(32, 13).Value = -1929.242
# Row 102
$ws.Cells.Item(102, 8).Value = 2569.25
$ws.Cells.Item(102, 9).Value = 2292.3333
$ws.Cells.Item(102, 11).Value = 2292.3333
$ws.Cells.Item(102, 13).Value = -670.3332999999998
# Row 110
$ws.Cells.Item(110, 8).Value = 1844.0869
$ws.Cells.Item(110, 9).Value = 1930.7
$ws.Cells.Item(110, 11).Value = 1930.7
$ws.Cells.Item(110, 13).Value = 114.3

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 970.05554
$ws.Cells.Item(20, 10).Value = 1032.5555
$ws.Cells.Item(20, 12).Value = 1032.5555
$ws.Cells.Item(20, 14).Value = -1526.5555
# Row 64
$ws.Cells.Item(64, 8).Value = 1673.5
$ws.Cells.Item(64, 9).Value = 1433.3334
$ws.Cells.Item(64, 10).Value = 1817.6
$ws.Cells.Item(64, 11).Value = 1433.3334
$ws.Cells.Item(64, 12).Value = 1817.6
$ws.Cells.Item(64, 13).Value = -1208.3334
$ws.Cells.Item(64, 14).Value = -2267.6
# Row 67
$ws.Cells.Item(67, 8).Value = 1673.5
$ws.Cells.Item(67, 9).Value = 1433.3334
$ws.Cells.Item(67, 10).Value = 1817.6
$ws.Cells.Item(67, 11).Value = 1433.3334
$ws.Cells.Item(67, 12).Value = 1817.6
$ws.Cells.Item(67, 13).Value = -653.3334
$ws.Cells.Item(67, 14).Value = -3377.6
# Row 86
$ws.Cells.Item(86, 8).Value = 3236.3125
$ws.Cells.Item(86, 9).Value = 2992
$ws.Cells.Item(86, 10).Value = 3969.25
$ws.Cells.Item(86, 11).Value = 2992
$ws.Cells.Item(86, 12).Value = 3969.25
$ws.Cells.Item(86, 13).Value = -1869
$ws.Cells.Item(86, 14).Value = -6215.25
# Row 89
$ws.Cells.Item(89, 8).Value = 3236.3125
$ws.Cells.Item(89, 9).Value = 2992
$ws.Cells.Item(89, 10).Value = 3969.25
$ws.Cells.Item(89, 11).Value = 14960
$ws.Cells.Item(89, 12).Value = 19846.25
$ws.Cells.Item(89, 13).Value = -9344
$ws.Cells.Item(89, 14).Value = -31078.25
# Row 94
$ws.Cells.Item(94, 8).Value = 3217.1428
$ws.Cells.Item(94, 9).Value = 2731.889
$ws.Cells.Item(94, 10).Value = 4090.6
$ws.Cells.Item(94, 11).Value = 2731.889
$ws.Cells.Item(94, 12).Value = 4090.6
$ws.Cells.Item(94, 13).Value = -2280.889
$ws.Cells.Item(94, 14).Value = -4992.6
# Row 99
$ws.Cells.Item(99, 8).Value = 1447.1538
$ws.Cells.Item(99, 9).Value = 1164.8182
$ws.Cells.Item(99, 11).Value = 1164.8182
$ws.Cells.Item(99, 13).Value = 333.1818000000001
# Row 107
$ws.Cells.Item(107, 8).Value = 2014.1818
$ws.Cells.Item(107, 9).Value = 1965.6
$ws.Cells.Item(107, 11).Value = 1965.6
$ws.Cells.Item(107, 13).Value = -45.59999999999991

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 427.7143
$ws.Cells.Item(22, 9).Value = 415.66666
$ws.Cells.Item(22, 11).Value = 415.66666
$ws.Cells.Item(22, 13).Value = -65.66665999999998
# Row 134
$ws.Cells.Item(134, 8).Value = 2788
$ws.Cells.Item(134, 9).Value = 2584.2
$ws.Cells.Item(134, 10).Value = 3399.4
$ws.Cells.Item(134, 11).Value = 7752.599999999999
$ws.Cells.Item(134, 12).Value = 10198.2
$ws.Cells.Item(134, 13).Value = -5217.599999999999
$ws.Cells.Item(134, 14).Value = -15268.2

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Cells.Item(4, 8).Value = 2341673.5
$ws.Cells.Item(4, 9).Value = 2653763.2
$ws.Cells.Item(4, 10).Value = 999.5
$ws.Cells.Item(4, 11).Value = 7961289.600000001
$ws.Cells.Item(4, 12).Value = 2998.5
$ws.Cells.Item(4, 13).Value = -7961177.600000001
$ws.Cells.Item(4, 14).Value = -3222.5
# Row 113
$ws.Cells.Item(113, 8).Value = 697.625
$ws.Cells.Item(113, 9).Value = 596.5
$ws.Cells.Item(113, 10).Value = 798.75
$ws.Cells.Item(113, 11).Value = 1789.5
$ws.Cells.Item(113, 12).Value = 2396.25
$ws.Cells.Item(113, 13).Value = 380.5
$ws.Cells.Item(113, 14).Value = -6736.25
# Row 122
$ws.Cells.Item(122, 8).Value = 398.6
$ws.Cells.Item(122, 9).Value = 248.25
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 11).Value = 2234.25
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = 215.75
$ws.Cells.Item(122, 14).Value = -13900
# Row 138
$ws.Cells.Item(138, 8).Value = 5747.25
$ws.Cells.Item(138, 10).Value = 5600
$ws.Cells.Item(138, 12).Value = 16800
$ws.Cells.Item(138, 14).Value = -27080

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 16413.572
$ws.Cells.Item(102, 9).Value = 18315.834
$ws.Cells.Item(102, 10).Value = 5000
$ws.Cells.Item(102, 11).Value = 18315.834
$ws.Cells.Item(102, 12).Value = 5000
$ws.Cells.Item(102, 13).Value = -16693.834
$ws.Cells.Item(102, 14).Value = -8244
# Row 122
$ws.Cells.Item(122, 8).Value = 49175.953
$ws.Cells.Item(122, 9).Value = 1029.8182
$ws.Cells.Item(122, 11).Value = 3089.4546
$ws.Cells.Item(122, 13).Value = -639.4546
# Row 126
$ws.Cells.Item(126, 8).Value = 2432.2222
$ws.Cells.Item(126, 9).Value = 2444.3333
$ws.Cells.Item(126, 11).Value = 7332.999899999999
$ws.Cells.Item(126, 13).Value = -4862.999899999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 3607.125
$ws.Cells.Item(40, 9).Value = 2776.4
$ws.Cells.Item(40, 10).Value = 4991.6665
$ws.Cells.Item(40, 11).Value = 2776.4
$ws.Cells.Item(40, 12).Value = 4991.6665
$ws.Cells.Item(40, 13).Value = -2640.4
$ws.Cells.Item(40, 14).Value = -5263.6665
# Row 43
$ws.Cells.Item(43, 8).Value = 340933.12
$ws.Cells.Item(43, 9).Value = 6949.5
$ws.Cells.Item(43, 10).Value = 392315.22
$ws.Cells.Item(43, 11).Value = 6949.5
$ws.Cells.Item(43, 12).Value = 392315.22
$ws.Cells.Item(43, 13).Value = -6756.5
$ws.Cells.Item(43, 14).Value = -392701.22
# Row 140
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).Value = $null

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 24
$ws.Cells.Item(24, 8).Value = 4599.8
$ws.Cells.Item(24, 10).Value = 4599.8
$ws.Cells.Item(24, 12).Value = 4599.8
$ws.Cells.Item(24, 14).Value = -5059.8
# Row 40
$ws.Cells.Item(40, 8).Value = 38500
$ws.Cells.Item(40, 10).Value = 38500
$ws.Cells.Item(40, 12).Value = 38500
$ws.Cells.Item(40, 14).Value = -38798
# Row 107
$ws.Cells.Item(107, 8).Value = 776.3333
$ws.Cells.Item(107, 9).Value = 776.3333
$ws.Cells.Item(107, 11).Value = 2328.9999
$ws.Cells.Item(107, 13).Value = -408.9998999999998
# Row 122
$ws.Cells.Item(122, 8).Value = 1711.5834
$ws.Cells.Item(122, 9).Value = 1893.2222
$ws.Cells.Item(122, 10).Value = 1166.6666
$ws.Cells.Item(122, 11).Value = 5679.6666
$ws.Cells.Item(122, 12).Value = 3499.9998
$ws.Cells.Item(122, 13).Value = -3229.6666
$ws.Cells.Item(122, 14).Value = -8399.9998
# Row 126
$ws.Cells.Item(126, 8).Value = 2184.8572
$ws.Cells.Item(126, 9).Value = 1916
$ws.Cells.Item(126, 11).Value = 5748
$ws.Cells.Item(126, 13).Value = -3278
# Row 132
$ws.Cells.Item(132, 8).Value = 2922
$ws.Cells.Item(132, 9).Value = 1941.8462
$ws.Cells.Item(132, 11).Value = 5825.5386
$ws.Cells.Item(132, 13).Value = -3295.5386
